# Generate Report for Handoff
# Refresh the handoff timestamps for the b72fc213 file (row 6 in each sheet):
#   - Overview!G6  "Latest HO Xliff Generate Date" -> newest handoff time across languages
#   - zh-cn!H6      "Latest Handoff Datetime" for the zh-cn handoff
#   - de-de!H6      "Latest Handoff Datetime" for the de-de handoff

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-08-16 16:38:10"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H6").Value = "2016-08-16 16:37:58"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H6").Value = "2016-08-16 16:38:10"
